$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.886.34"
$ws.Range("E2").Value = '  +1.71%  '

$ws.Range("D3").Value = "'2.584.41"
$ws.Range("E3").Value = '  +0.69%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = "'521.36"
$ws.Range("E5").Value = '  +1.03%  '

$ws.Range("D6").Value = "'139.12"
$ws.Range("E6").Value = '  -2.34%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = "'0.564"
$ws.Range("E8").Value = '  -0.31%  '

$ws.Range("D9").Value = "'2.596.86"
$ws.Range("E9").Value = '  +0.71%  '

$ws.Range("D10").Value = "'6.57"
$ws.Range("E10").Value = '  +0.10%  '

$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("E12").Value = '  +2.01%  '

$ws.Range("E13").Value = '  +3.26%  '

$ws.Range("D14").Value = "'3.044.10"
$ws.Range("E14").Value = '  +0.90%  '

$ws.Range("D15").Value = "'58.985.75"
$ws.Range("E15").Value = '  +1.84%  '

$ws.Range("E16").Value = '  +0.70%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("E17").Value = '  -0.53%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = "'2.567.72"
$ws.Range("E18").Value = '  -3.12%  '

$ws.Range("D19").Value = "'337.16"
$ws.Range("E19").Value = '  -0.90%  '

$ws.Range("E20").Value = '  +0.26%  '

$ws.Range("E21").Value = '  -1.03%  '

$ws.Range("D22").Value = "'6.52"
$ws.Range("E22").Value = '  +3.40%  '

$ws.Range("D24").Value = "'66.08"
$ws.Range("E24").Value = '  +1.09%  '

$ws.Range("E25").Value = '  +1.11%  '

$ws.Range("E26").Value = '  +0.57%  '

$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("D28").Value = "'7.01"
$ws.Range("E28").Value = '  +0.76%  '

$ws.Range("E29").Value = '  +0.05%  '

$ws.Range("D30").Value = "'0.0₃0724"
$ws.Range("E30").Value = '  -2.42%  '

$ws.Range("E31").Value = '  -5.49%  '

$ws.Range("E32").Value = '  +0.27%  '

$ws.Range("D33").Value = "'18.64"
$ws.Range("E33").Value = '  -0.07%  '

$ws.Range("D34").Value = "'149.13"
$ws.Range("E34").Value = '  -0.47%  '

$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("E36").Value = '  -1.98%  '

$ws.Range("D37").Value = "'36.79"
$ws.Range("E37").Value = '  +2.02%  '

$ws.Range("D38").Value = "'1.46"
$ws.Range("E38").Value = '  +0.75%  '

$ws.Range("D39").Value = "'0.825"
$ws.Range("E39").Value = '  -0.64%  '

$ws.Range("D40").Value = "'0.810"
$ws.Range("E40").Value = '  -6.95%  '

$ws.Range("D41").Value = "'3.51"
$ws.Range("E41").Value = '  -0.26%  '

$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = '  +0.19%  '

$ws.Range("D43").Value = "'272.12"
$ws.Range("E43").Value = '  +0.89%  '

$ws.Range("D44").Value = "'10.74"
$ws.Range("E44").Value = '  +0.81%  '

$ws.Range("D45").Value = "'0.0953"
$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D46").Value = "'0.589"
$ws.Range("E46").Value = '  +0.65%  '

$ws.Range("E47").Value = '  -0.62%  '

$ws.Range("D48").Value = "'18.43"
$ws.Range("E48").Value = '  -1.42%  '

$ws.Range("D49").Value = "'1.967.74"
$ws.Range("E49").Value = '  -0.53%  '

$ws.Range("D50").Value = "'4.51"
$ws.Range("E50").Value = '  -0.76%  '

$ws.Range("E51").Value = '  -0.24%  '
